$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Row 10 ("R40" rule): update the "From" value (C10) from 18 to 1
$ws.Range("C10").Value = 1
